$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.173228025436401
$ws.Range("B1").Value = 2.438494205474854
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.365235805511475
$ws.Range("E1").Value = 1.235641241073608
